# Allow n/a for judge non-availability (#6530)
# Adds a new judge row (row 10) to the "Judge Non-Availability Dates" sheet:
#   B10 = "Lamphere, Doris"  (judge name, text)
#   C10 = 861                (VLJ #, numeric)
#   D10 = "N/A"               (date column now holds the literal text "N/A"
#                              instead of a date, demonstrating that the
#                              upload format now tolerates a non-date value)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The judge-name and date columns need an explicit text format so that the
# written values ("Lamphere, Doris" / "N/A") are stored as shared strings
# rather than being coerced/interpreted, matching how "N/A" must survive as
# literal text in the Date column.
$ws.Range("B10").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"

$ws.Range("B10").Value = "Lamphere, Doris"
$ws.Range("C10").Value = 861
$ws.Range("D10").Value = "N/A"
